$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nädal 5")

# Fill in row 9 with the new time log entry
$ws.Range("A9").Value = 4
$ws.Range("B9").Value = 43524
$ws.Range("C9").Value = 0.54166666666666663
$ws.Range("D9").Value = 0.74652777777777779
$ws.Range("E9").Value = 30
$ws.Range("F9").FormulaArray = "=(D9-C9)*24*60 - E9"
$ws.Range("G9").Value = "Proge."
$ws.Range("H9").Value = "MVC EF"

# Update selection to A10
$ws.Range("A10").Select()

$wb.Save()
